$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(6)
$tbl = $sh.Table
$cell = $tbl.Cell(3, 1)
$tr = $cell.Shape.TextFrame.TextRange
$c1 = $tr.Characters(16, 8)
Write-Host "c1=[$($c1.Text)]"
$c1.Text = "ZUSNAHME"
$full = $cell.Shape.TextFrame.TextRange.Text
Write-Host "FULL=[$full]"
